# Add data for 2022-05-28
# Updates the "through May 19" snapshot to "through May 20" and bumps the
# May to-date carjacking counts for several neighborhoods/years.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet title / header label: May 19 -> May 20 -----------------------
$ws.Name = "Through 2022-05-20"
$ws.Range("B1").Value = "May 2022 (through May 20)"

# --- Updated counts on existing cells ------------------------------------
$ws.Range("B2").Value = 5     # Englewood, May 2022
$ws.Range("B4").Value = 3     # Humboldt Park, May 2022
$ws.Range("L4").Value = 4     # Humboldt Park, May 2020
$ws.Range("V4").Value = 3     # Humboldt Park, May 2018
$ws.Range("G5").Value = 3     # Garfield Park, May 2021
$ws.Range("L5").Value = 4     # Garfield Park, May 2020
$ws.Range("G8").Value = 4     # South Shore, May 2021
$ws.Range("AA8").Value = 2    # South Shore, May 2017
$ws.Range("G9").Value = 2     # Loop, May 2021
$ws.Range("G10").Value = 3    # Belmont Cragin, May 2021
$ws.Range("B25").Value = 4    # Auburn Gresham, May 2022
$ws.Range("AA30").Value = 3   # West Loop, May 2017

# --- New counts on cells that were previously blank ----------------------
$ws.Range("L6").Value = 1     # Chicago Lawn, May 2020
$ws.Range("AA6").Value = 1    # Chicago Lawn, May 2017
$ws.Range("L32").Value = 1    # United Center, May 2020
$ws.Range("Q32").Value = 1    # United Center, May 2019
$ws.Range("V39").Value = 1    # New City, May 2018
$ws.Range("B41").Value = 1    # Morgan Park, May 2022
$ws.Range("G58").Value = 1    # Clearing, May 2021
$ws.Range("L61").Value = 1    # East Village, May 2020
